$wb = $excel.ActiveWorkbook
$nl = [char]10

foreach ($ws in $wb.Worksheets) {

    # --- Row 10 (brand new row) -------------------------------------------
    # Build it from the pristine row 9 (A9:H9) so it inherits the exact same
    # per-column style pattern row 10 needs (s6,s7,s7,s7,s7,s7,s8,s7), then
    # overwrite the cell values that differ.
    $ws.Range("A9:H9").Copy($ws.Range("A10:H10"))
    $ws.Rows.Item(10).RowHeight = 60
    $ws.Range("A10").Value = "T1" + $nl + "(17:30-19:30)"
    $ws.Range("B10").Value = "Ngô Văn I"
    $ws.Range("G10").Value = "Lớp: CL10" + $nl + "Môn: Kỹ năng mềm" + $nl + "Phòng: R104" + $nl + "(Lý thuyết)"

    # --- Propagate the "highlighted" style (style 8) onto the new cells ---
    # Do this before the old style-8 donor cells (D8, G9) get converted back
    # to plain cells below.
    $ws.Range("D8").Copy($ws.Range("G8"))
    $ws.Range("D8").Copy($ws.Range("H8"))
    $ws.Range("G9").Copy($ws.Range("F9"))

    # --- Convert the old highlighted cells back to plain empty cells ------
    $ws.Range("C8").Copy($ws.Range("D8"))
    $ws.Range("D8").Value = ""
    $ws.Range("H9").Copy($ws.Range("G9"))
    $ws.Range("G9").Value = ""

    # --- Final text values for rows 8 and 9 --------------------------------
    $ws.Range("A8").Value = "S2" + $nl + "(09:00-11:00)"
    $ws.Range("G8").Value = "Lớp: CL10" + $nl + "Môn: Tiếng Anh chuyên ngành" + $nl + "Phòng: R102" + $nl + "(Lý thuyết)"
    $ws.Range("H8").Value = "Lớp: CL05" + $nl + "Môn: Tiếng Anh chuyên ngành" + $nl + "Phòng: R105" + $nl + "(Lý thuyết)"

    $ws.Range("A9").Value = "C1" + $nl + "(13:00-15:00)"
    $ws.Range("F9").Value = "Lớp: CL05" + $nl + "Môn: Kỹ năng mềm" + $nl + "Phòng: R103" + $nl + "(Lý thuyết)"
}
